# The commit adds a new weekly price record for "Cebollín" at
# "Vega Monumental Concepción". The new record is inserted as row 18,
# pushing the existing rows 18:110 down to 19:111 (dimension grows from
# A1:R110 to A1:R111).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18, shifting rows 18:110 -> 19:111.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new data record.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 45063
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100112037
$ws.Range("G18").Value = "Cebollín"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 3500
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = 3750
$ws.Range("N18").Value = "`$/paquete 36 unidades"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 104
$ws.Range("Q18").Value = 36
$ws.Range("R18").Value = "Hortaliza"

Write-Host "Inserted new row 18; dimension now covers A1:R$($ws.UsedRange.Rows.Count)"
